$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview": add a new row for c97bb94f-b051-4e6c-950d-75a882bd2ec0.md
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rOverview = $rowOverview.Range.Row()

$wsOverview.Cells.Item($rOverview, 1).Value = "c97bb94f-b051-4e6c-950d-75a882bd2ec0.md"
$wsOverview.Cells.Item($rOverview, 2).Value = "e2e\c97bb94f-b051-4e6c-950d-75a882bd2ec0.md"
$wsOverview.Cells.Item($rOverview, 3).Value = ".md"
$wsOverview.Cells.Item($rOverview, 4).Value = "'"
$wsOverview.Cells.Item($rOverview, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($rOverview, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($rOverview, 7).Value = "2016-11-03 19:09:15"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rOverview, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3bad24629534d886f76335b6f63a11e840da48df/e2e/c97bb94f-b051-4e6c-950d-75a882bd2ec0.md", "", "", "e2e\c97bb94f-b051-4e6c-950d-75a882bd2ec0.md")
$wsOverview.Cells.Item($rOverview, 2).Font.Underline = 2
$wsOverview.Cells.Item($rOverview, 2).Font.Color = 13272797

# -----------------------------------------------------------------
# Sheet "zh-cn": add a new row for c97bb94f-b051-4e6c-950d-75a882bd2ec0.md
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$rZh = $rowZh.Range.Row()

$wsZh.Cells.Item($rZh, 1).Value = "c97bb94f-b051-4e6c-950d-75a882bd2ec0.md"
$wsZh.Cells.Item($rZh, 2).Value = ".md"
$wsZh.Cells.Item($rZh, 3).Value = "Ready for handoff"
$wsZh.Cells.Item($rZh, 4).Value = "e2e"
$wsZh.Cells.Item($rZh, 5).Value = "ht"
$wsZh.Cells.Item($rZh, 6).Value = "False"
$wsZh.Cells.Item($rZh, 7).Value = "c97bb94f-b051-4e6c-950d-75a882bd2ec0.3bad24629534d886f76335b6f63a11e840da48df.zh-cn.xlf"
$wsZh.Cells.Item($rZh, 8).Value = "2016-11-03 19:09:02"
$wsZh.Cells.Item($rZh, 9).Value = "'"
$wsZh.Cells.Item($rZh, 10).Value = "'"
$wsZh.Cells.Item($rZh, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item($rZh, 12).Value = "'"
$wsZh.Cells.Item($rZh, 13).Value = "True"
$wsZh.Cells.Item($rZh, 14).Value = "'"
$wsZh.Cells.Item($rZh, 15).Value = "False"
$wsZh.Cells.Item($rZh, 16).Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rZh, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3bad24629534d886f76335b6f63a11e840da48df/e2e/c97bb94f-b051-4e6c-950d-75a882bd2ec0.md", "", "", "c97bb94f-b051-4e6c-950d-75a882bd2ec0.md")
$wsZh.Cells.Item($rZh, 1).Font.Underline = 2
$wsZh.Cells.Item($rZh, 1).Font.Color = 13272797

# -----------------------------------------------------------------
# Sheet "de-de": add a new row for c97bb94f-b051-4e6c-950d-75a882bd2ec0.md
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$rDe = $rowDe.Range.Row()

$wsDe.Cells.Item($rDe, 1).Value = "c97bb94f-b051-4e6c-950d-75a882bd2ec0.md"
$wsDe.Cells.Item($rDe, 2).Value = ".md"
$wsDe.Cells.Item($rDe, 3).Value = "Ready for handoff"
$wsDe.Cells.Item($rDe, 4).Value = "e2e"
$wsDe.Cells.Item($rDe, 5).Value = "ht"
$wsDe.Cells.Item($rDe, 6).Value = "False"
$wsDe.Cells.Item($rDe, 7).Value = "c97bb94f-b051-4e6c-950d-75a882bd2ec0.3bad24629534d886f76335b6f63a11e840da48df.de-de.xlf"
$wsDe.Cells.Item($rDe, 8).Value = "2016-11-03 19:09:15"
$wsDe.Cells.Item($rDe, 9).Value = "'"
$wsDe.Cells.Item($rDe, 10).Value = "'"
$wsDe.Cells.Item($rDe, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item($rDe, 12).Value = "'"
$wsDe.Cells.Item($rDe, 13).Value = "True"
$wsDe.Cells.Item($rDe, 14).Value = "'"
$wsDe.Cells.Item($rDe, 15).Value = "False"
$wsDe.Cells.Item($rDe, 16).Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rDe, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3bad24629534d886f76335b6f63a11e840da48df/e2e/c97bb94f-b051-4e6c-950d-75a882bd2ec0.md", "", "", "c97bb94f-b051-4e6c-950d-75a882bd2ec0.md")
$wsDe.Cells.Item($rDe, 1).Font.Underline = 2
$wsDe.Cells.Item($rDe, 1).Font.Color = 13272797

Write-Host "Report generated for handoff"
